$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row restructuring: the row that carries the taller (ht=51) row height
#    moves from row 6 (action/gremlin row) to row 5 (the header row), since
#    the header row wraps text in the two new "Final ..." header cells.
#    Deleting row 5 shifts the ht=51 row up to become row 5, then inserting a
#    fresh blank row pushes it back down to row 6 - leaving row 5 with ht=51
#    and row 6 with the default (no explicit) height, matching the target.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(6).Insert()

# Row 5 now contains left-over cells from the old row 6 (action/gremlin
# cells) that were shifted up into this position. Clear them out before
# writing the real header content.
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("P5").ClearContents()
$ws.Range("Q5").ClearContents()
$ws.Range("R5").ClearContents()
$ws.Range("S5").ClearContents()

# ---------------------------------------------------------------------------
# 2. Header row (row 5)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Game #"
$ws.Range("B5").Value = "Round 1"
$ws.Range("E5").Value = "Round 2"
$ws.Range("H5").Value = "Round 3"
$ws.Range("K5").Value = "Round 4"
$ws.Range("N5").Value = "…."
$ws.Range("O5").Value = "Gremlin Header"

$ws.Range("T5").Value = "Final Capacity"
$ws.Range("T5").WrapText = $true

$ws.Range("U5").Value = "Final User Story Chance"
$ws.Range("U5").WrapText = $true

$ws.Range("V5").Value = "Normalized score"

# ---------------------------------------------------------------------------
# 3. Data row (row 6)
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "Clarify Product Vision"
$ws.Range("C6").Value = "Protected from Outside Distraction"

$ws.Range("D6").Value = "Protected from Outside Distraction"
$ws.Range("D6").Font.ThemeColor = 1

$ws.Range("E6").Value = "All Work is Done on Main or Trunk"
$ws.Range("F6").Value = "Unit Testing"
$ws.Range("G6").Value = "Remote Team Avatars"

$ws.Range("H6").Value = "action"
$ws.Range("I6").Value = "action"
$ws.Range("J6").Value = "action"
$ws.Range("K6").Value = "action"
$ws.Range("L6").Value = "action"
$ws.Range("M6").Value = "action"

$ws.Range("O6").Value = "Poor Quality"
$ws.Range("P6").Value = "gremlin"
$ws.Range("Q6").Value = "gremlin"
$ws.Range("R6").Value = "gremlin"
$ws.Range("S6").Value = "gremlin"

$ws.Range("T6").Style = "Normal"
$ws.Range("U6").Style = "Normal"
$ws.Range("T6").Value = 45
$ws.Range("U6").Value = 75
$ws.Range("V6").Value = 173

# ---------------------------------------------------------------------------
# 4. New row 8 - a single, empty, center-aligned cell at E8
# ---------------------------------------------------------------------------
$ws.Range("E8").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Selection / active cell
# ---------------------------------------------------------------------------
$ws.Range("A5").Select()

Write-Host "done"
